# A new daily price record for "Apio" (Feria Lagunitas de Puerto Montt) is
# inserted as row 46. All existing records from row 46 down to row 194
# shift down by one row (to rows 47..195), and the sheet's used range
# grows from A1:R194 to A1:R195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46:194 down by one row, freeing up row 46 for the new record.
$ws.Rows("46:46").Insert()

# Populate the newly freed row 46 with the new record.
$ws.Range("A46").Value = 4
$ws.Range("B46").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C46").Value = "Los Lagos"
$ws.Range("D46").Value = 44592
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = 100112017
$ws.Range("G46").Value = "Apio"
$ws.Range("H46").Value = "Americana (o)"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 25
$ws.Range("K46").Value = 12000
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = 12000
$ws.Range("N46").Value = "`$/docena de matas"
$ws.Range("O46").Value = "Región de Coquimbo"
$ws.Range("P46").Value = 2000
$ws.Range("Q46").Value = 6
$ws.Range("R46").Value = "Hortaliza"
